$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 19.64514939341828
$ws.Range("C2").Value = 10.24460485648591
$ws.Range("D2").Value = 13.74373863765531
$ws.Range("E2").Value = 14.35017290282688
$ws.Range("G2").Value = 37.46021275869425
$ws.Range("H2").Value = 16.54923131936297
$ws.Range("J2").Value = 8.677308883221299
$ws.Range("L2").Value = 12.14041658435898
$ws.Range("O2").Value = 26.33518898919076
$ws.Range("B3").Value = 19.08561092465424
$ws.Range("C3").Value = 10.02182573314117
$ws.Range("D3").Value = 13.75975593837957
$ws.Range("E3").Value = 14.39948570684341
$ws.Range("G3").Value = 37.63089685077181
$ws.Range("H3").Value = 16.62015612213224
$ws.Range("J3").Value = 8.692301114567117
$ws.Range("L3").Value = 12.11181271498297
$ws.Range("O3").Value = 26.45852599297091
$ws.Range("B4").Value = 18.73492722965919
$ws.Range("C4").Value = 9.882442860901492
$ws.Range("D4").Value = 13.77226126301739
$ws.Range("E4").Value = 14.43189931566359
$ws.Range("G4").Value = 37.74883244059805
$ws.Range("H4").Value = 16.66689202831964
$ws.Range("J4").Value = 8.701990485781526
$ws.Range("L4").Value = 12.0956507070449
$ws.Range("O4").Value = 26.54083349730091
$ws.Range("B5").Value = 18.59042664907353
$ws.Range("C5").Value = 9.825053724840114
$ws.Range("D5").Value = 13.77802865145838
$ws.Range("E5").Value = 14.44564547665045
$ws.Range("G5").Value = 37.80017634692021
$ws.Range("H5").Value = 16.68673842536982
$ws.Range("J5").Value = 8.706061084804469
$ws.Range("L5").Value = 12.08942084314307
$ws.Range("O5").Value = 26.57602429998184
$ws.Range("B6").Value = 18.56634219099742
$ws.Range("C6").Value = 9.815490603911083
$ws.Range("D6").Value = 13.7790268616362
$ws.Range("E6").Value = 14.44796048292421
$ws.Range("G6").Value = 37.80889976382884
$ws.Range("H6").Value = 16.69008227805296
$ws.Range("J6").Value = 8.706744390356432
$ws.Range("L6").Value = 12.08840801055087
$ws.Range("O6").Value = 26.58196722586556
$ws.Range("B7").Value = 18.73298462804736
$ws.Range("C7").Value = 9.881671192064829
$ws.Range("D7").Value = 13.77233632607963
$ws.Range("E7").Value = 14.43208252493989
$ws.Range("G7").Value = 37.74951161004361
$ws.Range("H7").Value = 16.66715644061358
$ws.Range("J7").Value = 8.702044888421259
$ws.Range("L7").Value = 12.09556524119908
$ws.Range("O7").Value = 26.54130141839423
$ws.Range("B8").Value = 19.45381353524598
$ws.Range("C8").Value = 10.16836630656851
$ws.Range("D8").Value = 13.74870705414048
$ws.Range("E8").Value = 14.36673300789719
$ws.Range("G8").Value = 37.51632743968788
$ws.Range("H8").Value = 16.57302428722697
$ws.Range("J8").Value = 8.682377960221688
$ws.Range("L8").Value = 12.13026569738345
$ws.Range("O8").Value = 26.37634773717608
$ws.Range("B9").Value = 20.80285055138988
$ws.Range("C9").Value = 10.70736225256412
$ws.Range("D9").Value = 13.72356224326857
$ws.Range("E9").Value = 14.2555076830973
$ws.Range("G9").Value = 37.16408067327695
$ws.Range("H9").Value = 16.41374571945458
$ws.Range("J9").Value = 8.647634956352567
$ws.Range("L9").Value = 12.20924324806745
$ws.Range("O9").Value = 26.10525889086875
$ws.Range("B10").Value = 21.74498446359927
$ws.Range("C10").Value = 11.085943769972
$ws.Range("D10").Value = 13.71800039334186
$ws.Range("E10").Value = 14.1840831689365
$ws.Range("G10").Value = 36.97037283495112
$ws.Range("H10").Value = 16.31218436075741
$ws.Range("J10").Value = 8.624416246347565
$ws.Range("L10").Value = 12.27367068023479
$ws.Range("O10").Value = 25.93828448079653
$ws.Range("B11").Value = 22.16118129321432
$ws.Range("C11").Value = 11.2537666663491
$ws.Range("D11").Value = 13.71826811017877
$ws.Range("E11").Value = 14.15382013500643
$ws.Range("G11").Value = 36.8965949062996
$ws.Range("H11").Value = 16.26934516553696
$ws.Range("J11").Value = 8.614349262098651
$ws.Range("L11").Value = 12.30431055565544
$ws.Range("O11").Value = 25.86936954283701
$ws.Range("B12").Value = 22.31687468345459
$ws.Range("C12").Value = 11.31663768643736
$ws.Range("D12").Value = 13.71877091192426
$ws.Range("E12").Value = 14.14268040185957
$ws.Range("G12").Value = 36.87073538858224
$ws.Range("H12").Value = 16.25360713062296
$ws.Range("J12").Value = 8.610607996778322
$ws.Range("L12").Value = 12.31609907598566
$ws.Range("O12").Value = 25.84429056438745
$ws.Range("B13").Value = 22.28343027056663
$ws.Range("C13").Value = 11.30312824684521
$ws.Range("D13").Value = 13.71864478797492
$ws.Range("E13").Value = 14.14506530404648
$ws.Range("G13").Value = 36.87621200197979
$ws.Range("H13").Value = 16.25697504765317
$ws.Range("J13").Value = 8.611410598102605
$ws.Range("L13").Value = 12.31355203223565
$ws.Range("O13").Value = 25.84964643489706
$ws.Range("B14").Value = 22.17402928840955
$ws.Range("C14").Value = 11.25895298934386
$ws.Range("D14").Value = 13.71830143684394
$ws.Range("E14").Value = 14.15289724592284
$ws.Range("G14").Value = 36.89442567979606
$ws.Range("H14").Value = 16.26804067813691
$ws.Range("J14").Value = 8.614040047304758
$ws.Range("L14").Value = 12.30527671214163
$ws.Range("O14").Value = 25.86728585778027
$ws.Range("B15").Value = 22.10676547060519
$ws.Range("C15").Value = 11.23180447370946
$ws.Range("D15").Value = 13.71814337082954
$ws.Range("E15").Value = 14.15773623260247
$ws.Range("G15").Value = 36.90585324659079
$ws.Range("H15").Value = 16.27488178401528
$ws.Range("J15").Value = 8.615659881785639
$ws.Range("L15").Value = 12.30023187486196
$ws.Range("O15").Value = 25.87822318163423
$ws.Range("B16").Value = 21.71752436236208
$ws.Range("C16").Value = 11.07488352961407
$ws.Range("D16").Value = 13.7180391251143
$ws.Range("E16").Value = 14.18610575092582
$ws.Range("G16").Value = 36.97548439220324
$ws.Range("H16").Value = 16.31505171841404
$ws.Range("J16").Value = 8.625084085166353
$ws.Range("L16").Value = 12.27169461126566
$ws.Range("O16").Value = 25.94293036825323
$ws.Range("B17").Value = 21.47546855694605
$ws.Range("C17").Value = 10.9774566224149
$ws.Range("D17").Value = 13.71869119680055
$ws.Range("E17").Value = 14.20408010849104
$ws.Range("G17").Value = 37.02188582128703
$ws.Range("H17").Value = 16.34055624886387
$ws.Range("J17").Value = 8.630992147412625
$ws.Range("L17").Value = 12.25452503397896
$ws.Range("O17").Value = 25.9844333994874
$ws.Range("B18").Value = 21.33508570017282
$ws.Range("C18").Value = 10.92100811959621
$ws.Range("D18").Value = 13.71932966468248
$ws.Range("E18").Value = 14.21462823703927
$ws.Range("G18").Value = 37.04992359109264
$ws.Range("H18").Value = 16.35554213551416
$ws.Range("J18").Value = 8.634436951104053
$ws.Range("L18").Value = 12.24477522435545
$ws.Range("O18").Value = 26.00896743334939
$ws.Range("B19").Value = 21.28735964915674
$ws.Range("C19").Value = 10.90182653059144
$ws.Range("D19").Value = 13.71959110415786
$ws.Range("E19").Value = 14.21823568408432
$ws.Range("G19").Value = 37.0596478570404
$ws.Range("H19").Value = 16.36067040811354
$ws.Range("J19").Value = 8.635611323869975
$ws.Range("L19").Value = 12.24149585893947
$ws.Range("O19").Value = 26.01738788584661
$ws.Range("B20").Value = 21.50135672433132
$ws.Range("C20").Value = 10.98787079015619
$ws.Range("D20").Value = 13.71859452615068
$ws.Range("E20").Value = 14.20214500030327
$ws.Range("G20").Value = 37.01680658248399
$ws.Range("H20").Value = 16.33780850037764
$ws.Range("J20").Value = 8.630358399268763
$ws.Range("L20").Value = 12.25633979764238
$ws.Range("O20").Value = 25.97994671992346
$ws.Range("B21").Value = 22.2062158560511
$ws.Range("C21").Value = 11.2719471474899
$ws.Range("D21").Value = 13.71839140120993
$ws.Range("E21").Value = 14.15058812567032
$ws.Range("G21").Value = 36.88901934303441
$ws.Range("H21").Value = 16.26477728526622
$ws.Range("J21").Value = 8.613265793313193
$ws.Range("L21").Value = 12.30770237561558
$ws.Range("O21").Value = 25.86207707511137
$ws.Range("B22").Value = 22.65569000193764
$ws.Range("C22").Value = 11.45362505462575
$ws.Range("D22").Value = 13.72059786770492
$ws.Range("E22").Value = 14.11875912997595
$ws.Range("G22").Value = 36.81762548128211
$ws.Range("H22").Value = 16.21987008621254
$ws.Range("J22").Value = 8.602507790823498
$ws.Range("L22").Value = 12.34235140991281
$ws.Range("O22").Value = 25.79097611881959
$ws.Range("B23").Value = 22.41686100950623
$ws.Range("C23").Value = 11.35703937725257
$ws.Range("D23").Value = 13.7192065535836
$ws.Range("E23").Value = 14.13557615518059
$ws.Range("G23").Value = 36.85461535371161
$ws.Range("H23").Value = 16.24357932021831
$ws.Range("J23").Value = 8.608211862225119
$ws.Range("L23").Value = 12.32376159478416
$ws.Range("O23").Value = 25.82837949768467
$ws.Range("B24").Value = 21.48965647929188
$ws.Range("C24").Value = 10.98316390080398
$ws.Range("D24").Value = 13.71863740985682
$ws.Range("E24").Value = 14.20301919476104
$ws.Range("G24").Value = 37.01909866821021
$ws.Range("H24").Value = 16.33904975140119
$ws.Range("J24").Value = 8.630644766697943
$ws.Range("L24").Value = 12.25551896501921
$ws.Range("O24").Value = 25.98197305051712
$ws.Range("B25").Value = 20.44586077649695
$ws.Range("C25").Value = 10.5643855578076
$ws.Range("D25").Value = 13.72809572058108
$ws.Range("E25").Value = 14.28378779296351
$ws.Range("G25").Value = 37.24801651996806
$ws.Range("H25").Value = 16.45412188410289
$ws.Range("J25").Value = 8.656626999283921
$ws.Range("L25").Value = 12.18673356560071
$ws.Range("O25").Value = 26.17296059128045
